$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.14 = 11919.01 pesos`n✅ 11919.01 pesos = 3.12 = 965.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 318.399
$ws2.Range("O10").Value = 3795
$ws2.Range("N12").Value = 3822.5
$ws2.Range("O12").Value = 309.58
